# Apply the "upto 6 nov" commit to the "locators" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("locators")

# --- Write every new/changed cell value, in the exact order the original
# author entered them (this keeps the shared-string table append order
# identical to the authored workbook). Formatting is applied afterwards,
# since it does not touch the shared-string table. ---
$ws.Cells.Item(14, 1).Value2 = "add_to_bag_sun_glass_item"
$ws.Cells.Item(15, 1).Value2 = "cart_added_pop_up"
$ws.Cells.Item(15, 2).Value2 = '(By.XPATH, "//div[@class=''minicart-prod-name'']")'
$ws.Cells.Item(16, 1).Value2 = "ajio_logo_home"
$ws.Cells.Item(16, 2).Value2 = '(By.XPATH, "//img[@alt=''Ajio logo'']")'
$ws.Cells.Item(17, 1).Value2 = "cart_button_home_page"
$ws.Cells.Item(14, 2).Value2 = '(By.XPATH, "//span[text()=''ADD TO BAG'']")'
$ws.Cells.Item(17, 2).Value2 = '(By.XPATH, "//a/div[@class=''ic-cart '']")'
$ws.Cells.Item(10, 2).Value2 = '(By.LINK_TEXT, "Sunglasses & Frames")'
$ws.Cells.Item(12, 1).Value2 = "sun_glass_item_head"
$ws.Cells.Item(12, 2).Value2 = '(By.XPATH, "//h1[@class=''prod-name'']")'
$ws.Cells.Item(13, 1).Value2 = "sun_glass_item_price"
$ws.Cells.Item(13, 2).Value2 = '(By.XPATH, "//div[@class=''prod-sp'']")'
$ws.Cells.Item(20, 1).Value2 = "proceed_to_bag_btn"
$ws.Cells.Item(20, 2).Value2 = '(By.XPATH, "//div[@aria-label=''PROCEED TO BAG'']")'
$ws.Cells.Item(18, 1).Value2 = "cart_quantity_fetch"
$ws.Cells.Item(18, 2).Value2 = '(By.XPATH, "(//div[@class=''minicart-value''])[3]")'
$ws.Cells.Item(19, 1).Value2 = "cart_price_fetch"
$ws.Cells.Item(19, 2).Value2 = '(By.XPATH, "//div[@class=''minicart-totalamt'']/span")'
$ws.Cells.Item(21, 1).Value2 = "delete_btn_ship_page"
$ws.Cells.Item(21, 2).Value2 = '(By.XPATH, "//a/span[text()=''ROYAL SON'']/../../../..//div[@class=''product-delete'']/div[contains(text(), ''Delete'')]")'
$ws.Cells.Item(22, 1).Value2 = "item_price_shipping_page"
$ws.Cells.Item(22, 2).Value2 = '(By.XPATH, "//a/span[text()=''ROYAL SON'']/../../../..//div[@class=''priceinfo'']/div[2]")'

# --- Formatting for the new rows (12-22): column B keeps the sheet's
# 12pt body font, rows get the 15.5pt height used throughout the table,
# and four rows (the newest "sun glass item" + "cart fetch" locators)
# get a yellow highlight fill. ---
$highlightRows = @(12, 13, 18, 19)
for ($r = 12; $r -le 22; $r++) {
    $ws.Cells.Item($r, 2).Font.Size = 12
    $ws.Rows.Item($r).RowHeight = 15.5
    if ($highlightRows -contains $r) {
        $ws.Cells.Item($r, 1).Interior.Color = 65535
        $ws.Cells.Item($r, 2).Interior.Color = 65535
    }
}

# --- Column B width grew to fit the new, longer locator text ---
$ws.Columns.Item(2).ColumnWidth = 78.5

# --- Selection / view ---
$ws.Range("B21").Select() | Out-Null
